$d = $word.ActiveDocument

function Get-ParagraphContainingText($searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $matchStart = $r.Start
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $matchStart -and $matchStart -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

# -----------------------------------------------------------------
# 1) Remove the _GoBack bookmark from its current location (it sits
#    right after the "Asymmetry is there in policy..." paragraph).
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}

# -----------------------------------------------------------------
# 2) In the "Optimal policy is large..." paragraph, insert a new
#    sentence right before the trailing " (Jenny)", then re-create
#    the _GoBack bookmark right after the new sentence (i.e. right
#    before " (Jenny)").
# -----------------------------------------------------------------
$pJenny = Get-ParagraphContainingText("close to zero region.")

$rJenny = $pJenny.Range.Duplicate
$rJenny.Find.Execute(" (Jenny)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rJenny.InsertBefore(" It spends more than 60% of time in the <90 bp movements in int rate")

$rBookmark = $pJenny.Range.Duplicate
$rBookmark.Find.Execute(" (Jenny)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rBookmark.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rBookmark)

# -----------------------------------------------------------------
# 3) Add a new bullet paragraph right before "IRFs with cgain..."
# -----------------------------------------------------------------
$pIRF = Get-ParagraphContainingText("IRFs with cgain against unanchored")

$rBeforeIRF = $pIRF.Range.Duplicate
$rBeforeIRF.Collapse(1)
$rBeforeIRF.InsertBefore("Goodfriend 1993 showed that to subdue inflation scares, Fed repatedly raised bu hundreds of bp, e.g. raised by 300 bp in March 1980`r")

# -----------------------------------------------------------------
# 4) Add two new bullet paragraphs at the very end of the document,
#    right after "IRFs with cgain..." (which remains the last
#    paragraph since the previous insert went *before* it).
# -----------------------------------------------------------------
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertAfter("`rTR less aggressive than RE: refer to the debate between Orphanides vs. Preston and emphasize my contribution (Philippe)`rSame: Gurkaynak, Sack and Swanson 2005 find that interest rates in the future go negative after a positive int rate shock today, b/c corr(int, E(pi future) < 0).")
